$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 232.66667
$ws.Range("I33").Value = 227.85715
$ws.Range("K33").Value = 227.85715
$ws.Range("M33").Value = 1.14285000000001
$ws.Range("H70").Value = 9133.066000000001
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 9499.714
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 28499.142
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -29039.142
$ws.Range("H73").Value = 9133.066000000001
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 9499.714
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 28499.142
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -30371.142
$ws.Range("H104").Value = 1027.2
$ws.Range("I104").Value = 1027.2
$ws.Range("K104").Value = 3081.6
$ws.Range("M104").Value = -1334.6
$ws.Range("H132").Value = 3476
$ws.Range("I132").Value = 3258.2856
$ws.Range("K132").Value = 9774.856800000001
$ws.Range("M132").Value = -7244.856800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6266.857
$ws.Range("J2").Value = 16918.334
$ws.Range("L2").Value = 16918.334
$ws.Range("N2").Value = -17144.334
$ws.Range("H5").Value = 99
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H32").Value = 6627
$ws.Range("I32").Value = 4673.3687
$ws.Range("K32").Value = 4673.3687
$ws.Range("M32").Value = -4386.3687
$ws.Range("H97").Value = 1113.625
$ws.Range("I97").Value = 756.8570999999999
$ws.Range("K97").Value = 756.8570999999999
$ws.Range("M97").Value = -260.8570999999999
$ws.Range("H116").Value = 6266.857
$ws.Range("J116").Value = 16918.334
$ws.Range("L116").Value = 16918.334
$ws.Range("N116").Value = -21506.334
$ws.Range("H123").Value = 49999.5
$ws.Range("J123").Value = 49999.5
$ws.Range("L123").Value = 49999.5
$ws.Range("N123").Value = -59799.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6266.857
$ws.Range("J3").Value = 16918.334
$ws.Range("L3").Value = 16918.334
$ws.Range("N3").Value = -17146.334
$ws.Range("H4").Value = 99
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H86").Value = 4359
$ws.Range("I86").Value = 1791.5
$ws.Range("K86").Value = 1791.5
$ws.Range("M86").Value = -668.5
$ws.Range("H89").Value = 4359
$ws.Range("I89").Value = 1791.5
$ws.Range("K89").Value = 8957.5
$ws.Range("M89").Value = -3341.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 82.78570999999999
$ws.Range("I7").Value = 82.7
$ws.Range("J7").Value = 83
$ws.Range("K7").Value = 82.7
$ws.Range("L7").Value = 83
$ws.Range("M7").Value = 30.3
$ws.Range("N7").Value = -309
$ws.Range("H59").Value = 27109.285
$ws.Range("I59").Value = 6333.3335
$ws.Range("J59").Value = 42691.25
$ws.Range("K59").Value = 6333.3335
$ws.Range("L59").Value = 42691.25
$ws.Range("M59").Value = -5188.3335
$ws.Range("N59").Value = -44981.25
$ws.Range("H82").Value = 55000
$ws.Range("J82").Value = 55000
$ws.Range("L82").Value = 55000
$ws.Range("N82").Value = -55722
$ws.Range("H85").Value = 55000
$ws.Range("J85").Value = 55000
$ws.Range("L85").Value = 55000
$ws.Range("N85").Value = -57496
$ws.Range("H88").Value = 24833
$ws.Range("J88").Value = 24833
$ws.Range("L88").Value = 24833
$ws.Range("N88").Value = -25645
$ws.Range("H91").Value = 24833
$ws.Range("J91").Value = 24833
$ws.Range("L91").Value = 24833
$ws.Range("N91").Value = -27641

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 224.14285
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 228.16667
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 684.50001
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -3680.50001
$ws.Range("H130").Value = 3249
$ws.Range("I130").Value = 3249
$ws.Range("K130").Value = 9747
$ws.Range("M130").Value = -4727

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1411.6
$ws.Range("J80").Value = 3006
$ws.Range("L80").Value = 3006
$ws.Range("N80").Value = -5002
$ws.Range("H83").Value = 1411.6
$ws.Range("J83").Value = 3006
$ws.Range("L83").Value = 15030
$ws.Range("N83").Value = -25014
$ws.Range("H97").Value = 1229.091
$ws.Range("I97").Value = 1185
$ws.Range("J97").Value = 1282
$ws.Range("K97").Value = 1185
$ws.Range("L97").Value = 1282
$ws.Range("M97").Value = -689
$ws.Range("N97").Value = -2274

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4896.3125
$ws.Range("I7").Value = 3323.182
$ws.Range("K7").Value = 3323.182
$ws.Range("M7").Value = -3211.182
$ws.Range("H16").Value = 2994
$ws.Range("I16").Value = 659
$ws.Range("J16").Value = 9999
$ws.Range("K16").Value = 659
$ws.Range("L16").Value = 9999
$ws.Range("M16").Value = -489
$ws.Range("N16").Value = -10339
$ws.Range("H20").Value = 337001.34
$ws.Range("J20").Value = 500499.5
$ws.Range("L20").Value = 500499.5
$ws.Range("N20").Value = -500951.5
$ws.Range("H22").Value = 1017.06665
$ws.Range("I22").Value = 933.25
$ws.Range("J22").Value = 1112.8572
$ws.Range("K22").Value = 933.25
$ws.Range("L22").Value = 1112.8572
$ws.Range("M22").Value = -638.25
$ws.Range("N22").Value = -1702.8572
$ws.Range("H27").Value = 1017.06665
$ws.Range("I27").Value = 933.25
$ws.Range("J27").Value = 1112.8572
$ws.Range("K27").Value = 933.25
$ws.Range("L27").Value = 1112.8572
$ws.Range("M27").Value = -826.25
$ws.Range("N27").Value = -1326.8572
$ws.Range("H60").Value = 15500
$ws.Range("J60").Value = 15500
$ws.Range("L60").Value = 15500
$ws.Range("N60").Value = -16518
$ws.Range("H68").Value = 7079.4
$ws.Range("I68").Value = 5132.6665
$ws.Range("J68").Value = 9999.5
$ws.Range("K68").Value = 5132.6665
$ws.Range("L68").Value = 9999.5
$ws.Range("M68").Value = -4383.6665
$ws.Range("N68").Value = -11497.5
$ws.Range("H71").Value = 7079.4
$ws.Range("I71").Value = 5132.6665
$ws.Range("J71").Value = 9999.5
$ws.Range("K71").Value = 25663.3325
$ws.Range("L71").Value = 49997.5
$ws.Range("M71").Value = -21919.3325
$ws.Range("N71").Value = -57485.5
$ws.Range("H126").Value = 4896.3125
$ws.Range("I126").Value = 3323.182
$ws.Range("K126").Value = 9969.545999999998
$ws.Range("M126").Value = -7499.545999999998
$ws.Range("H130").Value = 22196.8
$ws.Range("J130").Value = 22196.8
$ws.Range("L130").Value = 22196.8
$ws.Range("N130").Value = -32236.8

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H59").Value = 32999
$ws.Range("J59").Value = 32999
$ws.Range("L59").Value = 32999
$ws.Range("N59").Value = -34475
$ws.Range("H81").Value = 819
$ws.Range("J81").Value = 699
$ws.Range("L81").Value = 1398
$ws.Range("N81").Value = -3520
$ws.Range("H84").Value = 819
$ws.Range("J84").Value = 699
$ws.Range("L84").Value = 6990
$ws.Range("N84").Value = -17598
$ws.Range("H136").Value = 5477.3335
$ws.Range("I136").Value = 5243.4614
$ws.Range("J136").Value = 6997.5
$ws.Range("K136").Value = 15730.3842
$ws.Range("L136").Value = 20992.5
$ws.Range("M136").Value = -13180.3842
$ws.Range("N136").Value = -26092.5
